$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Strip the "Heading2" paragraph style from the five section headings.
#    The runtime's Range.WordOpenXML is read-only, so InsertXML (which
#    REPLACES the range's own contents with the supplied OOXML) is used to
#    rewrite each heading paragraph with a bare <w:p> (no <w:pPr>/<w:pStyle>),
#    preserving the paragraph's existing run text exactly.
# ---------------------------------------------------------------------------
function Remove-HeadingStyle($para) {
    $t = $para.Range.Text
    $t = $t -replace "[\x07\x0d]+$", ""
    $escaped = $t -replace '&', '&amp;' -replace '<', '&lt;' -replace '>', '&gt;'
    $spaceAttr = ""
    if ($t -match "^\s" -or $t -match "\s$") {
        $spaceAttr = ' xml:space="preserve"'
    }
    $xml = '<?xml version="1.0" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:body><w:p><w:r><w:t' + $spaceAttr + '>' + $escaped + '</w:t></w:r></w:p></w:body></w:document>' +
        '</pkg:xmlData></pkg:part></pkg:package>'
    $para.Range.InsertXML($xml)
}

Remove-HeadingStyle($d.Paragraphs.Item(2))   # Introduction
Remove-HeadingStyle($d.Paragraphs.Item(6))   # Systemic Barriers
Remove-HeadingStyle($d.Paragraphs.Item(12))  # Impact on Students
Remove-HeadingStyle($d.Paragraphs.Item(18))  # Initiatives and Solutions
Remove-HeadingStyle($d.Paragraphs.Item(22))  # Conclusion

# ---------------------------------------------------------------------------
# 2) Swap the in-text citations for their updated sources. Replacements are
#    scoped to individual paragraphs (via Paragraph.Range.Find) because the
#    same old citation label maps to different new labels depending on
#    which paragraph it appears in.
# ---------------------------------------------------------------------------
function Replace-InParagraph($para, [string]$old, [string]$new) {
    $rng = $para.Range
    $rng.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

$enDash = [char]0x2013

# Paragraph 8: "Systemic barriers such as ..." (Systemic Barriers section)
Replace-InParagraph $d.Paragraphs.Item(8) "Pierszalowski et al." ("Nguyen 59" + $enDash + "60")
Replace-InParagraph $d.Paragraphs.Item(8) "Whitcomb and Singh" ("Nguyen 59" + $enDash + "60")

# Paragraph 10: "Additionally, educational disparities ..."
Replace-InParagraph $d.Paragraphs.Item(10) "Park et al." "Ref-u010943"
Replace-InParagraph $d.Paragraphs.Item(10) "Dickens" "Ref-u010943"

# Paragraph 14: "The challenges faced by ..." (Impact on Students section)
Replace-InParagraph $d.Paragraphs.Item(14) "Whitcomb and Singh" "Johnson"

# Paragraph 16: "Furthermore, disparities in STEM degree completion rates ..."
Replace-InParagraph $d.Paragraphs.Item(16) "Park et al." "Nguyen, 2015"

# Paragraph 20: "Initiatives such as mentorship programs ..."
Replace-InParagraph $d.Paragraphs.Item(20) "Dickens" "Ref-s577235"
Replace-InParagraph $d.Paragraphs.Item(20) "Sickle et al." "Ref-s577235"
